$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text format so numeric-looking strings
# (e.g. "66.433.28", "1.00", "7.05") are preserved verbatim as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '66.433.28'
$ws.Range("E2").Value = '  +2.84%  '
$ws.Range("D3").Value = '3.241.54'
$ws.Range("E3").Value = '  +5.19%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '578.68'
$ws.Range("E5").Value = '  +2.68%  '
$ws.Range("D6").Value = '155.48'
$ws.Range("E6").Value = '  +8.80%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.228.47'
$ws.Range("E8").Value = '  +5.21%  '
$ws.Range("E9").Value = '  +4.35%  '
$ws.Range("D10").Value = '7.05'
$ws.Range("E10").Value = '  +10.16%  '
$ws.Range("D11").Value = '0.166'
$ws.Range("E11").Value = '  +5.23%  '
$ws.Range("E12").Value = '  +4.43%  '
$ws.Range("D13").Value = '37.89'
$ws.Range("E13").Value = '  +5.65%  '
$ws.Range("E14").Value = '  +4.34%  '
$ws.Range("D15").Value = '3.756.30'
$ws.Range("E15").Value = '  +5.01%  '
$ws.Range("D16").Value = '561.47'
$ws.Range("E16").Value = '  +13.29%  '
$ws.Range("D17").Value = '66.449.77'
$ws.Range("E17").Value = '  +2.89%  '
$ws.Range("D18").Value = '3.239.01'
$ws.Range("E18").Value = '  +4.95%  '
$ws.Range("E19").Value = '  +3.17%  '
$ws.Range("E20").Value = '  +6.22%  '
$ws.Range("D21").Value = '14.44'
$ws.Range("E21").Value = '  +4.60%  '
$ws.Range("D22").Value = '0.745'
$ws.Range("E22").Value = '  +7.65%  '
$ws.Range("D23").Value = '7.88'
$ws.Range("E23").Value = '  +8.53%  '
$ws.Range("E24").Value = '  +7.03%  '
$ws.Range("D25").Value = '82.22'
$ws.Range("E25").Value = '  +3.93%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '9.45'
$ws.Range("E27").Value = '  +17.84%  '
$ws.Range("E28").Value = '  +6.20%  '
$ws.Range("E29").Value = '  +7.72%  '
$ws.Range("D30").Value = '27.97'
$ws.Range("E30").Value = '  +5.68%  '
$ws.Range("E31").Value = '  +2.54%  '
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("D34").Value = '565.38'
$ws.Range("E34").Value = '  +8.98%  '
$ws.Range("D35").Value = '5.79'
$ws.Range("E35").Value = '  +4.69%  '
$ws.Range("D36").Value = '6.42'
$ws.Range("E36").Value = '  +7.23%  '
$ws.Range("D37").Value = '0.0462'
$ws.Range("E37").Value = '  +13.68%  '
$ws.Range("D38").Value = '55.97'
$ws.Range("E38").Value = '  +4.53%  '
$ws.Range("D39").Value = '0.0871'
$ws.Range("E39").Value = '  +8.27%  '
$ws.Range("E40").Value = '  +14.00%  '
$ws.Range("E41").Value = '  +5.71%  '
$ws.Range("D42").Value = '3.139.77'
$ws.Range("E42").Value = '  +6.75%  '
$ws.Range("E43").Value = '  +3.01%  '
$ws.Range("D44").Value = '0.276'
$ws.Range("E44").Value = '  +11.15%  '
$ws.Range("E45").Value = '  +7.60%  '
$ws.Range("D46").Value = '26.70'
$ws.Range("E46").Value = '  +5.04%  '
$ws.Range("D47").Value = '0.0₃0563'
$ws.Range("E47").Value = '  +3.78%  '
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("E49").Value = '  +3.96%  '
$ws.Range("D50").Value = '2.27'
$ws.Range("E50").Value = '  +9.13%  '
$ws.Range("D51").Value = '122.57'
$ws.Range("E51").Value = '  +1.97%  '

# Restore the original (default) cell style now that values are set as text.
$ws.Range("D2:E51").Style = "Normal"
